# Update "Generate Report for Handback" timestamps.
# Workbook has 3 sheets: Overview, zh-cn, de-de.
# Row 3 on every sheet corresponds to file 55b57630-96f8-4b35-86f4-4d813023a6b7
# These date/time values are stored as plain text strings (shared strings, t="s"),
# so assigning plain strings keeps them as text and preserves the existing
# cell styling/number format.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: G3 = Latest HO Xliff Generate Date
$wsOverview.Range("G3").Value = "2016-08-23 10:47:22"

# zh-cn sheet: H3 = Correspond Handoff Datetime, K3 = Correspond Handback DateTime
$wsZhCn.Range("H3").Value = "2016-08-23 10:47:17"
$wsZhCn.Range("K3").Value = "2016-08-23 10:47:35"

# de-de sheet: H3 = Correspond Handoff Datetime, K3 = Correspond Handback DateTime
$wsDeDe.Range("H3").Value = "2016-08-23 10:47:22"
$wsDeDe.Range("K3").Value = "2016-08-23 10:47:42"
